# Auto-generated edit script: refresh cryptos price table values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "51.096.28"
$ws.Cells.Item(2, 5).Value = "  -1.15%  "
# Row 3
$ws.Cells.Item(3, 4).Value = "2.910.33"
$ws.Cells.Item(3, 5).Value = "  -0.63%  "
# Row 4
$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 5).Value = "  -0.10%  "
# Row 5
$ws.Cells.Item(5, 4).Value = "'371.10"
$ws.Cells.Item(5, 5).Value = "  +5.39%  "
# Row 6
$ws.Cells.Item(6, 4).Value = "'103.13"
$ws.Cells.Item(6, 5).Value = "  -3.55%  "
# Row 7
$ws.Cells.Item(7, 4).Value = "'0.539"
$ws.Cells.Item(7, 5).Value = "  -2.87%  "
# Row 8
$ws.Cells.Item(8, 4).Value = "'1.00"
$ws.Cells.Item(8, 5).Value = "  +0.07%  "
# Row 9
$ws.Cells.Item(9, 4).Value = "'0.586"
$ws.Cells.Item(9, 5).Value = "  -4.05%  "
# Row 10
$ws.Cells.Item(10, 4).Value = "'36.58"
$ws.Cells.Item(10, 5).Value = "  -3.29%  "
# Row 11
$ws.Cells.Item(11, 5).Value = "  +0.95%  "
# Row 12
$ws.Cells.Item(12, 4).Value = "'0.0834"
$ws.Cells.Item(12, 5).Value = "  -2.25%  "
# Row 13
$ws.Cells.Item(13, 4).Value = "'18.34"
$ws.Cells.Item(13, 5).Value = "  -3.27%  "
# Row 14
$ws.Cells.Item(14, 4).Value = "3.372.99"
# Row 15
$ws.Cells.Item(15, 4).Value = "'7.36"
$ws.Cells.Item(15, 5).Value = "  -2.77%  "
# Row 16
$ws.Cells.Item(16, 4).Value = "2.921.58"
$ws.Cells.Item(16, 5).Value = "  -0.35%  "
# Row 17
$ws.Cells.Item(17, 4).Value = "'0.934"
$ws.Cells.Item(17, 5).Value = "  -4.14%  "
# Row 18
$ws.Cells.Item(18, 4).Value = "51.040.72"
$ws.Cells.Item(18, 5).Value = "  -1.19%  "
# Row 19
$ws.Cells.Item(19, 4).Value = "'3.24"
$ws.Cells.Item(19, 5).Value = "  -7.16%  "
# Row 20
$ws.Cells.Item(20, 4).Value = "'7.20"
$ws.Cells.Item(20, 5).Value = "  -2.34%  "
# Row 21
$ws.Cells.Item(21, 4).Value = "'12.85"
$ws.Cells.Item(21, 5).Value = "  -4.56%  "
# Row 22
$ws.Cells.Item(22, 4).Value = "0.0₃0942"
$ws.Cells.Item(22, 5).Value = "  -2.20%  "
# Row 23
$ws.Cells.Item(23, 4).Value = "'68.14"
$ws.Cells.Item(23, 5).Value = "  -1.30%  "
# Row 24
$ws.Cells.Item(24, 4).Value = "'259.33"
$ws.Cells.Item(24, 5).Value = "  -1.00%  "
# Row 25
$ws.Cells.Item(25, 5).Value = "  -0.87%  "
# Row 26
$ws.Cells.Item(26, 2).Value = "LEO"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(26, 4).Value = "'4.34"
$ws.Cells.Item(26, 5).Value = "  +4.05%  "
# Row 27
$ws.Cells.Item(27, 2).Value = "Kaspa"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(27, 4).Value = "'0.171"
$ws.Cells.Item(27, 5).Value = "  -0.96%  "
# Row 28
$ws.Cells.Item(28, 2).Value = "Dai"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(28, 4).Value = "'1.00"
$ws.Cells.Item(28, 5).Value = "  -0.04%  "
# Row 29
$ws.Cells.Item(29, 2).Value = "EthereumClassic"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(29, 4).Value = "'25.67"
$ws.Cells.Item(29, 5).Value = "  -3.32%  "
# Row 30
$ws.Cells.Item(30, 2).Value = "Filecoin"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(30, 4).Value = "'6.94"
$ws.Cells.Item(30, 5).Value = "  -7.84%  "
# Row 31
$ws.Cells.Item(31, 2).Value = "Hedera"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(31, 4).Value = "'0.103"
$ws.Cells.Item(31, 5).Value = "  -1.17%  "
# Row 32
$ws.Cells.Item(32, 2).Value = "RenderToken"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(32, 4).Value = "'6.13"
$ws.Cells.Item(32, 5).Value = "  -0.19%  "
# Row 33
$ws.Cells.Item(33, 2).Value = "Cosmos"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(33, 4).Value = "'9.88"
$ws.Cells.Item(33, 5).Value = "  -3.48%  "
# Row 34
$ws.Cells.Item(34, 2).Value = "Toncoin"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(34, 4).Value = "'2.13"
$ws.Cells.Item(34, 5).Value = "  -1.44%  "
# Row 35
$ws.Cells.Item(35, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(35, 4).Value = "'34.62"
$ws.Cells.Item(35, 5).Value = "  -2.72%  "
# Row 36
$ws.Cells.Item(36, 2).Value = "OKB"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(36, 4).Value = "'50.82"
$ws.Cells.Item(36, 5).Value = "  -0.22%  "
# Row 37
$ws.Cells.Item(37, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(37, 4).Value = "'1.00"
$ws.Cells.Item(37, 5).Value = "  +0.37%  "
# Row 38
$ws.Cells.Item(38, 2).Value = "VeChain"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(38, 4).Value = "'0.0421"
$ws.Cells.Item(38, 5).Value = "  -1.65%  "
# Row 39
$ws.Cells.Item(39, 2).Value = "LidoDAOToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(39, 4).Value = "'3.03"
$ws.Cells.Item(39, 5).Value = "  -3.52%  "
# Row 40
$ws.Cells.Item(40, 2).Value = "Stacks"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(40, 4).Value = "'2.65"
$ws.Cells.Item(40, 5).Value = "  -0.63%  "
# Row 41
$ws.Cells.Item(41, 2).Value = "Celestia"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(41, 4).Value = "'17.06"
$ws.Cells.Item(41, 5).Value = "  -3.48%  "
# Row 42
$ws.Cells.Item(42, 2).Value = "ARBITRUM"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(42, 4).Value = "'1.84"
$ws.Cells.Item(42, 5).Value = "  -6.35%  "
# Row 43
$ws.Cells.Item(43, 2).Value = "Stellar"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(43, 4).Value = "'0.113"
$ws.Cells.Item(43, 5).Value = "  -2.41%  "
# Row 44
$ws.Cells.Item(44, 2).Value = "EnergySwap"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(44, 4).Value = "'22.08"
$ws.Cells.Item(44, 5).Value = "  -2.50%  "
# Row 45
$ws.Cells.Item(45, 2).Value = "Monero"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(45, 4).Value = "'119.43"
$ws.Cells.Item(45, 5).Value = "  +0.13%  "
# Row 46
$ws.Cells.Item(46, 2).Value = "WEMIXToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(46, 4).Value = "'2.10"
$ws.Cells.Item(46, 5).Value = "  -2.30%  "
# Row 47
$ws.Cells.Item(47, 2).Value = "Maker"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(47, 4).Value = "2.019.65"
$ws.Cells.Item(47, 5).Value = "  -3.96%  "
# Row 48
$ws.Cells.Item(48, 2).Value = "ApeXProtocol"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(48, 4).Value = "'2.30"
$ws.Cells.Item(48, 5).Value = "  -1.51%  "
# Row 49
$ws.Cells.Item(49, 2).Value = "NEARProtocol"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(49, 4).Value = "'3.17"
$ws.Cells.Item(49, 5).Value = "  -4.59%  "
# Row 50
$ws.Cells.Item(50, 2).Value = "TheGraph"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(50, 4).Value = "'0.243"
$ws.Cells.Item(50, 5).Value = "  +1.31%  "
# Row 51
$ws.Cells.Item(51, 2).Value = "BEAM"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Cells.Item(51, 4).Value = "'0.0310"
$ws.Cells.Item(51, 5).Value = "  -9.93%  "
